$p = $ppt.ActivePresentation
$s = $p.Slides.Add(3, 12)
try {
$shp2 = $s.Shapes.AddPicture("/tmp/work/does_not_exist.png", $false, $true, 60,0,50,50)
Write-Host "added: " $shp2.Name
} catch {
Write-Host "ERR: $_"
}
